$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColToNum($col) {
    $n = 0
    foreach ($ch in $col.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

# ---------------------------------------------------------------------------
# 1. Report-generation timestamp (row 5)
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:48 AM"

# ---------------------------------------------------------------------------
# 2. Report-summary totals (rows 8-9) - updated to reflect the new week's data
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 18178.79
$ws.Range("C9").Value = 45

# ---------------------------------------------------------------------------
# 3. Monday section, Point 01 CON-10-AAA-1-B-REEL line (row 16) - quantity x4
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 444
$ws.Range("H16").Value = 1038.96

# ---------------------------------------------------------------------------
# 4. Monday section TOTAL (row 29) - reflects the row 16 change above
# ---------------------------------------------------------------------------
$ws.Range("H29").Value = 4200.849999999999

# ---------------------------------------------------------------------------
# 5. Append a whole new day section - "Wednesday (07/30/2025)" - rows 32-66.
#    First stamp down the correct cell styles by copying formats from the
#    analogous, already-styled Monday block (rows 14-29), then overwrite the
#    values with the Wednesday data.
# ---------------------------------------------------------------------------

# Row 32 banner <- format of row 14 banner
$ws.Range("A14:H14").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null

# Row 33 column headers <- format of row 15 column headers
$ws.Range("A15:H15").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null

# Rows 34-65 (32 data rows) <- repeating 2-row striped format of rows 16-17
$ws.Range("A16:H17").Copy() | Out-Null
$ws.Range("A34:H65").PasteSpecial(-4122) | Out-Null

# Row 66 TOTAL <- format of row 29 TOTAL
$ws.Range("A29:H29").Copy() | Out-Null
$ws.Range("A66").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Now populate the values/text for every new cell, row by row.
# ---------------------------------------------------------------------------
$newRows = @(
    @{Row=32; Cells=@(@{Col="A"; Val="Wednesday (07/30/2025)"})},
    @{Row=33; Cells=@(@{Col="A"; Val="Point Number"}, @{Col="B"; Val="Billable Unit Code"}, @{Col="C"; Val="Work Type"}, @{Col="D"; Val="Unit Description"}, @{Col="E"; Val="Unit of Measure"}, @{Col="F"; Val="# Units"}, @{Col="G"; Val="N/A"}, @{Col="H"; Val="Pricing"})},
    @{Row=34; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="ANC-DHM-10-84-T1-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="ANC,Dbl Hlx Mach,10in,84in,TpEye 1in,Cor"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=217.53})},
    @{Row=35; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="ARM-8SF-GN-DL-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="ARM,8ftSgl.Fiberglass,Gain,DE Light,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=350.53})},
    @{Row=36; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="ARM-8SF-GN-TL-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="ARM,8ft Sgl.Fiberglass,Gain,Tang LD,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=350.53})},
    @{Row=37; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="GYF-38-D-78P-EP-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="GYF,3/8,Down,78in Pole mt,EyePlate,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=79.34999999999999})},
    @{Row=38; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="INS-15-P-S-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="INS,15kV,Pin,Silicon Polymer,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=5}, @{Col="H"; Val=470.85})},
    @{Row=39; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="PIN-15-PTP-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="Pin,15kV,Pole top,Corrosive"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=94.17})},
    @{Row=40; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="PIN-XAL-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="Pin,Crossarm Light,Corrosive"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=5}, @{Col="H"; Val=470.85})},
    @{Row=41; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="POL-45-2"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="Pole,45ft,Class 2"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=478.55})},
    @{Row=42; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="SAA-3-CV-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="SAA,3 inch,Clevis,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=55.18})},
    @{Row=43; Cells=@(@{Col="A"; Val="Point 01"}, @{Col="B"; Val="SAA-DE-20-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="SAA,DE Clamp #4-2/0, Corr."}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=55.18})},
    @{Row=44; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="ANC-DHM-10-84-T1-C"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="ANC,Dbl Hlx Mach,10in,84in,TpEye 1in,Cor"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=28.6})},
    @{Row=45; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="ARM-8D-60S"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="ARM,8ft Double Wood,60in Sgl.Wd"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=73.73})},
    @{Row=46; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="ARM-8S-60S"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="ARM,8ft Single Wood,60in Sgl.Wd"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=61.83})},
    @{Row=47; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="CON-2-AAA-1-B-REEL"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="CON,#2 AWG,Alum Alloy,One,Bare,Reels"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=480}, @{Col="H"; Val=446.4})},
    @{Row=48; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="DEC-20AL-C"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="DEC,#4 - #2/0 AA,AL,AS,Corrosive"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=3}, @{Col="H"; Val=285.45})},
    @{Row=49; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="GYF-38-D-42W-GH-C"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="GYF,3/8,Down,42in Wire mt,Guy Hook,Corr."}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=35.58})},
    @{Row=50; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="INS-15-P-S-C"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="INS,15kV,Pin,Silicon Polymer,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=5}, @{Col="H"; Val=155.4})},
    @{Row=51; Cells=@(@{Col="A"; Val="Point 02"}, @{Col="B"; Val="PIN-XAL-C"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="Pin,Crossarm Light,Corrosive"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=5}, @{Col="H"; Val=155.4})},
    @{Row=52; Cells=@(@{Col="A"; Val="Point 03"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=426}, @{Col="H"; Val=996.84})},
    @{Row=53; Cells=@(@{Col="A"; Val="Point 05"}, @{Col="B"; Val="CON-10-AAA-1-B"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=115}, @{Col="H"; Val=269.1})},
    @{Row=54; Cells=@(@{Col="A"; Val="Point 05"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=789}, @{Col="H"; Val=1846.26})},
    @{Row=55; Cells=@(@{Col="A"; Val="Point 05"}, @{Col="B"; Val="SAA-3-CV-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="SAA,3 inch,Clevis,Corr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=2}, @{Col="H"; Val=110.36})},
    @{Row=56; Cells=@(@{Col="A"; Val="Point 15"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=53}, @{Col="H"; Val=124.02})},
    @{Row=57; Cells=@(@{Col="A"; Val="Point 16"}, @{Col="B"; Val="CON-2-AAA-1-B-REEL"}, @{Col="C"; Val="Rem"}, @{Col="D"; Val="CON,#2 AWG,Alum Alloy,One,Bare,Reels"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=108}, @{Col="H"; Val=100.44})},
    @{Row=58; Cells=@(@{Col="A"; Val="Point 17"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=433}, @{Col="H"; Val=1013.22})},
    @{Row=59; Cells=@(@{Col="A"; Val="Point 17"}, @{Col="B"; Val="EQL-3-4-C-40-S-X"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="EQL,3 Ph,#4,CU Sol,4/0,CU Str,Xfr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=670})},
    @{Row=60; Cells=@(@{Col="A"; Val="Point 19"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=539}, @{Col="H"; Val=1261.26})},
    @{Row=61; Cells=@(@{Col="A"; Val="Point 19"}, @{Col="B"; Val="EQL-3-4-C-50-S-X"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="EQL,3 Ph,#4,CU Sol,500,CU Str,Xfr"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=1}, @{Col="H"; Val=670})},
    @{Row=62; Cells=@(@{Col="A"; Val="Point 21"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=231}, @{Col="H"; Val=540.54})},
    @{Row=63; Cells=@(@{Col="A"; Val="Point 23"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=618}, @{Col="H"; Val=1446.12})},
    @{Row=64; Cells=@(@{Col="A"; Val="Point 25"}, @{Col="B"; Val="CON-10-AAA-1-B-REEL"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="CON,#1/0 AWG,Alum Alloy,One,Bare,Reel"}, @{Col="E"; Val="FT"}, @{Col="F"; Val=88}, @{Col="H"; Val=205.92})},
    @{Row=65; Cells=@(@{Col="A"; Val="Point 27"}, @{Col="B"; Val="DEC-20AL-C"}, @{Col="C"; Val="Inst"}, @{Col="D"; Val="DEC,#4 - #2/0 AA,AL,AS,Corrosive"}, @{Col="E"; Val="EA"}, @{Col="F"; Val=3}, @{Col="H"; Val=858.75})},
    @{Row=66; Cells=@(@{Col="A"; Val="TOTAL"}, @{Col="H"; Val=13977.94})},
)

foreach ($rowDef in $newRows) {
    $r = $rowDef.Row
    foreach ($cellDef in $rowDef.Cells) {
        $colNum = ColToNum $cellDef.Col
        $ws.Cells.Item($r, $colNum).Value = $cellDef.Val
    }
}

# ---------------------------------------------------------------------------
# 6. New merged cells for the Wednesday section banner + TOTAL row
# ---------------------------------------------------------------------------
$ws.Range("A32:H32").Merge() | Out-Null
$ws.Range("A66:G66").Merge() | Out-Null

Write-Host "Wednesday (07/30/2025) section added; summary totals synced."
